$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "b" exponent value (B6) from 6 to 4 -- this cascades through the
# dependent formulas (H1, H3, B7, E5, E7, H5, H7).
$ws.Range("B6").Value = 4

# Update the sheet view: scroll so column E is the left-most visible column,
# zoom to 120%, and select the full columns A:H (active cell H1).
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Application.ActiveWindow.Zoom = 120
$ws.Range("A1:H1048576").Select()
$ws.Application.ActiveCell = $ws.Range("H1")

# Column width changes: G loses "best fit" (fixed width 23), H widens to
# 32.88671875 (best fit keeps autofit on).
$ws.Columns.Item(7).ColumnWidth = 23
$ws.Columns.Item(8).ColumnWidth = 32.88671875

# Remove sheet protection.
$ws.Unprotect()
